# Washing Machine Design Calculations — "Updating for New Manufacturing Ideas"
#
# The applied normal force design input drops from 1000 lbf to 600 lbf.
# AppliedForceNormal is a defined name pointing at 'Design and Material
# Specs'!D5, so every downstream formula across the workbook (Motor Sizing
# Calcs, Piston Calcs, Frame Calcs, Screw Calculations, ...) recalculates
# from this single input.
$wb = $excel.ActiveWorkbook

$specs = $wb.Worksheets.Item("Design and Material Specs")
$specs.Range("D5").Value = 600

# Reflect the author's final on-screen state: "Screw Calculations" was the
# previously-active tab (selection parked on E25); "Motor Sizing Calcs"
# becomes the newly-active tab with its selection moved to E41. Touch the
# soon-to-be-inactive sheet's selection first, then activate + select on
# the sheet that ends up active so it is the last (and therefore winning)
# activation.
$screw = $wb.Worksheets.Item("Screw Calculations")
$screw.Activate()
$screw.Range("E25").Select()

$motor = $wb.Worksheets.Item("Motor Sizing Calcs")
$motor.Activate()
$motor.Range("E41").Select()
